# "Generate Report for Handoff"
#
# The b.md file has been handed off for localization (zh-cn and de-de).
# 1. Overview sheet: status for b.md flips from "Handed back: in sync
#    with en-US" to "Ready for handoff" for both locale columns.
# 2. zh-cn / de-de sheets: the "Latest Handoff File" (column C) and
#    "Latest Handoff Datetime" (column D) for the b.md row (row 3) are
#    updated to reflect the new handoff package + timestamp, and the
#    hyperlink backing column C now points at the new handoff file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-26 09:01:53"

$zhcnLink = $zhcn.Range("C3").Hyperlinks.Item(1)
$zhcnLink.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcnLink.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/041edb829a7b6d7c5e5b17a5e1d54fa97f69b19e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-01-26 09:02:04"

$dedeLink = $dede.Range("C3").Hyperlinks.Item(1)
$dedeLink.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dedeLink.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b44af1cb5171fd8b2cde915c25331db78e43d6c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

Write-Output "Report generated for handoff."
